$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the BF14 model-label: the old "model_13_3_1_..._20230622 ..."
#    shared string is retired and a new "model_13_3_1_..._20240430" string
#    takes its place. Style also reverts to the sheet default (no explicit
#    cell style), matching the author's re-entry of this label.
# ---------------------------------------------------------------------------
$ws.Range("BF14").ClearFormats()
$ws.Range("BF14").Value = "model_13_3_1_lr1e-3_ep300_2080Ti6_20240430"

# ---------------------------------------------------------------------------
# 2) Re-run/refresh the experiment numbers for the first model block
#    (rows 2-10, columns BF:BI). Row 2 additionally gets a new bold
#    Times New Roman look (fresh font + number formats), rows 3-10 keep
#    their existing look and only the values change.
# ---------------------------------------------------------------------------
$ws.Range("BF2").NumberFormat = "0.0000E+00"
$ws.Range("BF2").Font.Bold = $true
$ws.Range("BF2").Font.Name = "Times New Roman"
$ws.Range("BF2").Font.Color = 0

$ws.Range("BG2:BI2").NumberFormat = "0.0000"
$ws.Range("BG2:BI2").Font.Bold = $true
$ws.Range("BG2:BI2").Font.Name = "Times New Roman"
$ws.Range("BG2:BI2").Font.Color = 0

$ws.Range("BF2").Value = 0.0048453500494360898
$ws.Range("BG2").Value = 0.76571071147918701
$ws.Range("BH2").Value = 1.91854012012481
$ws.Range("BI2").Value = 4.0345716476440403

$ws.Range("BF3").Value = 0.0048944219015538597
$ws.Range("BG3").Value = 0.75839370489120395
$ws.Range("BH3").Value = 1.7165440320968599
$ws.Range("BI3").Value = 3.8926594257354701

$ws.Range("BF4").Value = 0.0051650004461407601
$ws.Range("BG4").Value = 0.76235479116439797
$ws.Range("BH4").Value = 1.7833321094512899
$ws.Range("BI4").Value = 3.9301307201385498

$ws.Range("BF5").Value = 0.0051110610365867597
$ws.Range("BG5").Value = 0.77026826143264704
$ws.Range("BH5").Value = 2.0108397006988499
$ws.Range("BI5").Value = 4.0277943611145002

$ws.Range("BF6").Value = 0.00454193167388439
$ws.Range("BG6").Value = 0.77186399698257402
$ws.Range("BH6").Value = 2.3355793952941801
$ws.Range("BI6").Value = 4.4158411026000897

$ws.Range("BF7").Value = 0.0045466427691280798
$ws.Range("BG7").Value = 0.78389561176300004
$ws.Range("BH7").Value = 2.73815441131591
$ws.Range("BI7").Value = 4.5856542587280202

$ws.Range("BF8").Value = 0.0073956199921667498
$ws.Range("BG8").Value = 0.62689977884292603
$ws.Range("BH8").Value = -1.8175415992736801
$ws.Range("BI8").Value = 2.1936058998107901

$ws.Range("BF9").Value = 0.0067749423906206998
$ws.Range("BG9").Value = 0.65076559782028198
$ws.Range("BH9").Value = -1.14425504207611
$ws.Range("BI9").Value = 2.4851877689361501

$ws.Range("BF10").Value = 0.0072098122909665099
$ws.Range("BG10").Value = 0.65893441438674905
$ws.Range("BH10").Value = -0.98290163278579701
$ws.Range("BI10").Value = 2.4970710277557302

# ---------------------------------------------------------------------------
# 3) Recalculate so the AVERAGE() roll-ups in rows 12/13 (BF:BI) pick up the
#    refreshed inputs automatically.
# ---------------------------------------------------------------------------
$excel.Calculate()

# ---------------------------------------------------------------------------
# 4) Window/view bookkeeping: the author scrolled/zoomed the sheet and left
#    the selection on the BH12:BI13 summary block.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Zoom = 115
$ws.Range("AX1").Select()
$ws.Range("BH12:BI13").Select()
